# ListadePrecioFormulario.xlsx - add Sucursales / Medicos / Compañias sections
# mirroring the existing "Estudios" block, plus the corresponding named ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

function Add-PriceBlock {
    param(
        [string]$TitleCell,
        [string]$TitleRange,
        [string]$Title,
        [string]$ClaveHeaderRange,
        [string]$ClaveHeaderCell,
        [string]$NombreHeaderRange,
        [string]$NombreHeaderCell,
        [string]$ListaHeaderRange,
        [string]$ListaHeaderCell,
        [string]$ClaveDataCell,
        [string]$NombreDataCell,
        [string]$ListaDataCell,
        [string]$BlankRange1,
        [string]$BlankRange2,
        [string]$BlankRange3
    )

    # Row 13 - section title (merged across the whole block)
    $r = $ws.Range($TitleRange)
    $r.HorizontalAlignment = $xlCenter
    $ws.Range($TitleCell).Value = $Title
    $r.Merge()

    # Row 14 - column headers (Clave / Nombre / Lista de Precios)
    $r = $ws.Range($ClaveHeaderRange)
    $r.HorizontalAlignment = $xlCenter
    $ws.Range($ClaveHeaderCell).Value = "Clave"
    $r.Merge()

    $r = $ws.Range($NombreHeaderRange)
    $r.HorizontalAlignment = $xlCenter
    $ws.Range($NombreHeaderCell).Value = "Nombre"
    $r.Merge()

    $r = $ws.Range($ListaHeaderRange)
    $r.HorizontalAlignment = $xlCenter
    $ws.Range($ListaHeaderCell).Value = "Lista de Precios"
    $r.Merge()

    # Row 15 - placeholder tokens for the template engine
    $ws.Range($ClaveDataCell).Value = "{{item.Clave}}"
    $ws.Range($NombreDataCell).Value = "{{item.Nombre}}"
    $ws.Range($ListaDataCell).Value = "{{item.ListaPrecio}}"

    # Row 16 - blank placeholder row (merged, centered, empty)
    $r = $ws.Range($BlankRange1)
    $r.HorizontalAlignment = $xlCenter
    $r.Merge()

    $r = $ws.Range($BlankRange2)
    $r.HorizontalAlignment = $xlCenter
    $r.Merge()

    $r = $ws.Range($BlankRange3)
    $r.HorizontalAlignment = $xlCenter
    $r.Merge()
}

# ---- Sucursales (columns J:O) ----
Add-PriceBlock -TitleCell "J13" -TitleRange "J13:O13" -Title "Sucursales" `
    -ClaveHeaderRange "J14:K14" -ClaveHeaderCell "J14" `
    -NombreHeaderRange "L14:M14" -NombreHeaderCell "L14" `
    -ListaHeaderRange "N14:O14" -ListaHeaderCell "N14" `
    -ClaveDataCell "J15" -NombreDataCell "L15" -ListaDataCell "N15" `
    -BlankRange1 "J16:K16" -BlankRange2 "L16:M16" -BlankRange3 "N16:O16"

# ---- Medicos (columns T:Y) ----
Add-PriceBlock -TitleCell "T13" -TitleRange "T13:Y13" -Title "Médicos" `
    -ClaveHeaderRange "T14:U14" -ClaveHeaderCell "T14" `
    -NombreHeaderRange "V14:W14" -NombreHeaderCell "V14" `
    -ListaHeaderRange "X14:Y14" -ListaHeaderCell "X14" `
    -ClaveDataCell "T15" -NombreDataCell "V15" -ListaDataCell "X15" `
    -BlankRange1 "T16:U16" -BlankRange2 "V16:W16" -BlankRange3 "X16:Y16"

# ---- Compañias (columns AB:AG) ----
Add-PriceBlock -TitleCell "AB13" -TitleRange "AB13:AG13" -Title "Compañias" `
    -ClaveHeaderRange "AB14:AC14" -ClaveHeaderCell "AB14" `
    -NombreHeaderRange "AD14:AE14" -NombreHeaderCell "AD14" `
    -ListaHeaderRange "AF14:AG14" -ListaHeaderCell "AF14" `
    -ClaveDataCell "AB15" -NombreDataCell "AD15" -ListaDataCell "AF15" `
    -BlankRange1 "AB16:AC16" -BlankRange2 "AD16:AE16" -BlankRange3 "AF16:AG16"

# Named ranges for the three new blocks (kept alongside the existing "Estudios" name)
$wb.Names.Add("Sucursales", "=Precios!`$J`$15:`$O`$16")
$wb.Names.Add("Medicos", "=Precios!`$T`$15:`$Y`$16")
$wb.Names.Add("Compañias", "=Precios!`$AB`$15:`$AG`$16")

# Leftover formatted (empty) cell + final selection, matching the saved workbook state
$ws.Range("O25").Font.Underline = $true
$ws.Range("O25").Select()

Write-Host "Done"
